$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: remove the "reordered to complete / group features first" comment in I24
# and let the row height return to the default (it was a custom 30pt height because
# of that wrapped comment).
$ws.Range("I24").ClearContents()
$ws.Rows.Item(24).AutoFit()

# Row 26: the "Front-end: Threejs scene - MQTT client" task was removed from the plan
# (merged back into the grouped trigger work), so blank out the whole row, keeping the
# date formatting on E26/F26 the same way the other blank rows (28, 29) do.
$ws.Range("A26:F26").ClearContents()

# Reflect the current selection/scroll position used while editing.
$ws.Range("F21").Select()
